# Chronos/testdata.xlsx edit
# - Rename "Skills Enhancement" -> "Skill Enhancement" in column D (Project)
# - Add a new "Task" column (H) with per-row task labels
# - Adjust column H width, selection, etc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Skills Enhancement" -> "Skill Enhancement" typo on every row that uses it ---
$ws.Range("D2").Value = "Skill Enhancement"
$ws.Range("D4").Value = "Skill Enhancement"
$ws.Range("D5").Value = "Skill Enhancement"
$ws.Range("D6").Value = "Skill Enhancement"
$ws.Range("D7").Value = "Skill Enhancement"

# --- Add the new "Task" column (H) ---
$ws.Range("H1").Value = "Task"
$ws.Range("H2").Value = "UFT Training"
$ws.Range("H3").Value = "Mandatory Fun"
$ws.Range("H4").Value = "UFT Training"
$ws.Range("H5").Value = "UFT Training"
$ws.Range("H6").Value = "UFT Training"
$ws.Range("H7").Value = "UFT Training"

# Match the bold header formatting used by the rest of row 1
$ws.Range("H1").Font.Bold = $true

# Widen the new column similarly to the others (closest attainable value)
$ws.Columns.Item(8).ColumnWidth = 21.7

# Move/update the active selection like the author's last click in the sheet
$ws.Range("H9").Select()
